$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting "une" (old row 5) and all the numbers below it down by one.
$ws.Rows.Item(5).Insert()

# New row 5 gets the "minuit" label, matching the same text style as the other label rows.
$ws.Range("A5").Value = "minuit"

# Reset the view so the first row is visible again, keeping B7 as the active cell.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select()
